$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dutybot")

# Insert a new row at the top of the sheet, shifting all existing rows down.
$ws.Rows.Item(1).Insert()

# Populate the new header row with "day" and " name".
$ws.Cells.Item(1, 1).Value = "day"
$ws.Cells.Item(1, 2).Value = " name"

# Match the formatting used by the data rows (now starting at row 2).
$ws.Range("A2:B2").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
